$d = $word.ActiveDocument

function Replace-ParagraphContent($paraIndex, $xml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = ""
    $rng.Collapse(1)
    $rng.InsertXML($xml)
}

# 1) "Implementation" heading: add lastRenderedPageBreak
Replace-ParagraphContent 30 '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Implementation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 2) Remove the _GoBack bookmark that sits after the flow-chart page break
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 3) "Security Protocol" heading: add lastRenderedPageBreak before "S"
Replace-ParagraphContent 34 '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:i/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:lastRenderedPageBreak/><w:t>S</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>ecurity</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> Protocol</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 4) "Code executions" heading: add lastRenderedPageBreak
Replace-ParagraphContent 39 '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Code executions</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 5) Split "Jetbrains" / "FileTransferApplication" runs with proofErr tags
Replace-ParagraphContent 42 '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">To run program via IDE (this was developed via IntelliJ IDEA from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Jetbrains</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">), first open </w:t></w:r><w:r><w:t xml:space="preserve">project </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>FileTransferApplication</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Then, individually run Server module and Client module. From this step just simply follow the prompt.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 6) Split "src" / "javac" runs with proofErr tags
Replace-ParagraphContent 43 '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>To run program</w:t></w:r><w:r><w:t xml:space="preserve"> in</w:t></w:r><w:r><w:t xml:space="preserve"> the command line environment</w:t></w:r><w:r><w:t>, first navigate to where the Server folder and the Client folder located on your local machine.</w:t></w:r><w:r><w:t xml:space="preserve"> In each module, navigate to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder. Run command &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Server.java&#8221; for Server and &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Client.java&#8221; for Client to compile the code. Then run &#8220;java Server&#8221; to start the Server and &#8220;java Client&#8221; to run the Client. After this, run as program instructs. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 7) Insert the large new "Detail about security protocol" / "Threat model" block
#    right after paragraph 43 and before the trailing empty paragraph.
$p43 = $d.Paragraphs.Item(43)
$insertionPoint = $d.Range($p43.Range.End, $p43.Range.End)
$insertionPoint.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
        </w:rPr>
        <w:t>Detail</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
        </w:rPr>
        <w:t>bout</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
        </w:rPr>
        <w:t>security</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
        </w:rPr>
        <w:t>protocol</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>uthentication</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Server</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>has</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>a</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> certificate which is its public key signed by CA’s private key. The client has CA’s public key. When server receiving the connection request from client, server will send the certificate to client. Client uses CA’s public key to verify the certificate and get server’s public key.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>fter authentication, client generates a random number. We consider this random number as a master key.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Encrypt master key with server’s public key and send it to server.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>C</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>onfidentiality</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">In this project, we use Sha1 to encrypt our files which is </w:t>
      </w:r>
      <w:r>
        <w:t>generating and using one-time pass to XOR with data.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> =</w:t>
      </w:r>
      <w:r>
        <w:t>Sha1(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>K</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>ab</w:t>
      </w:r>
      <w:r>
        <w:t>|IV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)                            c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>=p</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="微软雅黑" w:hint="eastAsia"/>
          <w:color w:val="333333"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>⊕</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> =Sha1(K</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>ab</w:t>
      </w:r>
      <w:r>
        <w:t>|c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:t>)                             c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:t>=p</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="微软雅黑" w:hint="eastAsia"/>
          <w:color w:val="333333"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>⊕</w:t>
      </w:r>
      <w:r>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">          …                                              …</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">i </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>=Sha1(</w:t>
      </w:r>
      <w:r>
        <w:t>K</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>ab</w:t>
      </w:r>
      <w:r>
        <w:t>|c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>i-1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>)                             c</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>i</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>=</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>p</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>i</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="微软雅黑" w:eastAsia="微软雅黑" w:hAnsi="微软雅黑" w:cs="微软雅黑" w:hint="eastAsia"/>
          <w:color w:val="333333"/>
          <w:sz w:val="20"/>
          <w:szCs w:val="20"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>⊕</w:t>
      </w:r>
      <w:r>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>i</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>ntegrity</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">We also use Sha1 to calculate the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>MAC(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>Message authentication codes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>MAC =SHA1(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>K</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="subscript"/>
        </w:rPr>
        <w:t>ab</w:t>
      </w:r>
      <w:r>
        <w:t>|m</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
        </w:rPr>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>hreat model</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>E</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>avesdrop</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Even if attacker can watch the message between client and server pass over the network, the attacker still can’t learn the contents of message between client and server.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>nitiate a conversation pretending to be Client</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>The project has one-way a</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>uthentication</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t xml:space="preserve">, so it can’t verify the client. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>L</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>ie in wait at Server’s network address and accept a connection from Client</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>The server will use certificate signed by CA to verify itself.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>R</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>ead Server’s database</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>The attacker can get server’s private key to decrypt the master key from client.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>Man-in-the-Middle Attack</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>The att</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>acker can’t get useful information if he doesn’t know server’s private key or master key.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:ind w:left="1440"/>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 8) Append the tail paragraphs (blank + blank + restored _GoBack bookmark)
#    after the final trailing empty paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$tailPoint.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:eastAsia="zh-CN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p/>
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# 9) Footer page number field cache: "1" -> "6"
$ftr = $d.Sections.Item(1).Footers.Item(1)
$ftr.Range.Characters.Item(1).Text = "6"

Write-Output "done"
